# Increase compatibility with India data + Tanzania pilot data
#
# 1. Insert a new "district" row right after the "a1-fid"/"hf_id" row (new row 9),
#    pushing the existing rows 9-38 down to 10-39.
# 2. Append a new "FormVersion" row at the end (new row 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new row 9 (district) ----------------------------------------
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "a1-district"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "district"
$ws.Range("D9").Value = "district"

# Match the formatting used by the other rows of this same kind:
#  - A9/C9 look like the other "a1-*" rows (e.g. row 8, "a1-fid")
#  - D9 looks like the other "raw == clean" rows that have distinct
#    highlighting (e.g. row 22, the "o1-o1_2a"/"o1_2a" row)
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D22").Copy()
$ws.Range("D9").PasteSpecial(-4122)

# --- 2. Append new row 40 (FormVersion) -------------------------------------
$ws.Range("A40").Value = "FormVersion"
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "form_version"
$ws.Range("D40").Value = "FormVersion"

# Match the formatting used by the other "simple" metadata rows (e.g. row 2,
# "date"/"date_call"/"date") for A/C. D40 naturally inherits the column's
# default look, which already matches.
$ws.Range("A2").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C40").PasteSpecial(-4122)

$excel.CutCopyMode = 0
